$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '43.215.49'
Set-TextValue 'E2' '  +0.72%  '
Set-TextValue 'D3' '2.342.41'
Set-TextValue 'E3' '  +1.36%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '302.63'
Set-TextValue 'E6' '  -1.57%  '
Set-TextValue 'D7' '0.505'
Set-TextValue 'E7' '  -0.13%  '
Set-TextValue 'E8' '  +0.03%  '
Set-TextValue 'E9' '  -0.66%  '
Set-TextValue 'D10' '34.16'
Set-TextValue 'E10' '  -2.27%  '
Set-TextValue 'E11' '  -0.71%  '
Set-TextValue 'D12' '18.71'
Set-TextValue 'E12' '  -3.41%  '
Set-TextValue 'E13' '  +2.26%  '
Set-TextValue 'D14' '6.76'
Set-TextValue 'E14' '  -1.26%  '
Set-TextValue 'D15' '2.706.46'
Set-TextValue 'E15' '  +1.39%  '
Set-TextValue 'D16' '2.368.54'
Set-TextValue 'E16' '  +3.08%  '
Set-TextValue 'D17' '0.799'
Set-TextValue 'E17' '  +1.68%  '
Set-TextValue 'D18' '43.139.18'
Set-TextValue 'E18' '  +0.67%  '
Set-TextValue 'D19' '12.19'
Set-TextValue 'E19' '  -2.34%  '
Set-TextValue 'D20' '6.21'
Set-TextValue 'E20' '  +3.13%  '
Set-TextValue 'D21' '0.0₃0892'
Set-TextValue 'E21' '  +0.04%  '
Set-TextValue 'D22' '68.04'
Set-TextValue 'E22' '  +0.54%  '
Set-TextValue 'D23' '236.01'
Set-TextValue 'E23' '  +0.17%  '
Set-TextValue 'D24' '2.22'
Set-TextValue 'E24' '  -0.74%  '
Set-TextValue 'E25' '  +0.01%  '
Set-TextValue 'E26' '  -0.34%  '
Set-TextValue 'D27' '24.65'
Set-TextValue 'E27' '  -0.19%  '
Set-TextValue 'D28' '2.23'
Set-TextValue 'E28' '  -6.01%  '
Set-TextValue 'D29' '9.19'
Set-TextValue 'E29' '  +1.31%  '
Set-TextValue 'D30' '31.56'
Set-TextValue 'E30' '  -3.11%  '
Set-TextValue 'D31' '1.00'
Set-TextValue 'E31' '  +0.02%  '
Set-TextValue 'E32' '  +0.82%  '
Set-TextValue 'D33' '0.0729'
Set-TextValue 'E33' '  +4.41%  '
Set-TextValue 'D34' '17.33'
Set-TextValue 'E34' '  -2.26%  '
Set-TextValue 'B35' 'ARBITRUM'
Set-TextValue 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '1.83'
Set-TextValue 'E35' '  +4.03%  '
Set-TextValue 'B36' 'RenderToken'
Set-TextValue 'C36' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D36' '4.38'
Set-TextValue 'E36' '  -2.41%  '
Set-TextValue 'E37' '  -0.87%  '
Set-TextValue 'D38' '0.101'
Set-TextValue 'E38' '  +0.52%  '
Set-TextValue 'B39' 'LidoDAOToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D39' '2.75'
Set-TextValue 'E39' '  -0.04%  '
Set-TextValue 'B40' 'EnergySwap'
Set-TextValue 'C40' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D40' '22.27'
Set-TextValue 'E40' '  +17.91%  '
Set-TextValue 'E41' '  -0.37%  '
Set-TextValue 'D42' '111.79'
Set-TextValue 'E42' '  -31.91%  '
Set-TextValue 'D43' '1.939.54'
Set-TextValue 'E43' '  -1.87%  '
Set-TextValue 'E44' '  +1.23%  '
Set-TextValue 'D45' '10.05'
Set-TextValue 'E45' '  -4.89%  '
Set-TextValue 'E46' '  +1.63%  '
Set-TextValue 'E47' '  -1.56%  '
Set-TextValue 'D48' '2.571.30'
Set-TextValue 'E48' '  +1.36%  '
Set-TextValue 'D49' '53.14'
Set-TextValue 'E49' '  -0.42%  '
Set-TextValue 'D50' '2.81'
Set-TextValue 'E50' '  -3.16%  '
Set-TextValue 'D51' '72.08'
Set-TextValue 'E51' '  -0.35%  '
